$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.619.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.118.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.17"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.117.26"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.69%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.648.78"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.62"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.729.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.117.02"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.67"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "360.95"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.36%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.06"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0865"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.29"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.10"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.48"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.68"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.511.91"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.64%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.697"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "37.80"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0270"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.975"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.73"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0912"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.64%  "
